$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "notes" header column (J1)
$ws.Range("J1").Value = "notes"

# New data rows from the May 2018 Prince Island trip
$ws.Range("A9").Value  = 2018
$ws.Range("B9").Value  = 5
$ws.Range("C9").Value  = 20180514
$ws.Range("D9").Value  = 1
$ws.Range("E9").Value  = 1
$ws.Range("F9").Value  = 21
$ws.Range("G9").Value  = 25
$ws.Range("H9").Value  = 22
$ws.Range("I9").Formula = "=AVERAGE(F9:H9)"

$ws.Range("A10").Value = 2018
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 20180514
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 73
$ws.Range("G10").Value = 75
$ws.Range("H10").Value = 72
$ws.Range("I10").Formula = "=AVERAGE(F10:H10)"

# Write the note strings in the order they first appeared upstream so the
# shared-string table indices line up (J10 "lower colony..." before J9
# "upper colony...").
$ws.Range("J10").Value = "lower colony on northern promentary of island, 3 more COMU right above colony"
$ws.Range("J9").Value  = "upper colony on northern promentary of island, 3 more COMU right above colony"

$ws.Range("A11").Value = 2018
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 20180514
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 32
$ws.Range("G11").Value = 32
$ws.Range("H11").Value = 33
$ws.Range("I11").Formula = "=AVERAGE(F11:H11)"

# Match the saved selection state from the author's session
$ws.Range("A10:C11").Select() | Out-Null

$wb.Save() | Out-Null
